$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the data range as Text first so numeric-looking strings (e.g.
# "103.00", "0.000007475", "6.310") are written verbatim instead of being
# coerced into numbers by Excels usual type inference.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "30.438.94"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "1.909.64"
$ws.Range("E3").Value = "  -1.53%  "
$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "0.9982"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "239.77"
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "0.9983"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "0.4767"
$ws.Range("E7").Value = "  -2.91%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.2846"
$ws.Range("E8").Value = "  -3.11%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.06701"
$ws.Range("E9").Value = "  -2.76%  "
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "19.41"
$ws.Range("E10").Value = "  +0.95%  "
$ws.Range("B11").Value = "Litecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D11").Value = "103.00"
$ws.Range("E11").Value = "  -2.20%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "0.07763"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.913.10"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "5.191"
$ws.Range("E14").Value = "  -3.12%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.6702"
$ws.Range("E15").Value = "  -4.69%  "
$ws.Range("B16").Value = "BitcoinCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D16").Value = "274.42"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "30.465.71"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "0.9982"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.000007475"
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "12.61"
$ws.Range("E20").Value = "  -3.59%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "5.374"
$ws.Range("E21").Value = "  -3.49%  "
$ws.Range("B22").Value = "BitDAO"
$ws.Range("C22").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D22").Value = "0.4649"
$ws.Range("E22").Value = "  -7.28%  "
$ws.Range("B23").Value = "BinanceUSD"
$ws.Range("C23").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D23").Value = "0.9977"
$ws.Range("E23").Value = "  -0.44%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "6.296"
$ws.Range("E24").Value = "  -3.62%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D25").Value = "9.329"
$ws.Range("E25").Value = "  -5.08%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "167.73"
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "19.19"
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.085"
$ws.Range("E28").Value = "  -2.90%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "1.379"
$ws.Range("E29").Value = "  -1.04%  "
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "0.09975"
$ws.Range("E30").Value = "  -4.25%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "4.583"
$ws.Range("E31").Value = "  +0.41%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "1.513"
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "4.268"
$ws.Range("E33").Value = "  -2.55%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.04716"
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.7268"
$ws.Range("E35").Value = "  -3.82%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").Value = "1.115"
$ws.Range("E36").Value = "  -3.13%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "2.714"
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.01904"
$ws.Range("E38").Value = "  -5.11%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "2.611"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "6.310"
$ws.Range("E40").Value = "  -2.41%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "74.07"
$ws.Range("E41").Value = "  -5.32%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "107.10"
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "1.957"
$ws.Range("E43").Value = "  -6.18%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.8604"
$ws.Range("E44").Value = "  -5.69%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").Value = "0.4265"
$ws.Range("E45").Value = "  -3.77%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "0.9975"
$ws.Range("E46").Value = "  -0.12%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.404"
$ws.Range("E47").Value = "  -3.30%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "955.84"
$ws.Range("E48").Value = "  -4.49%  "
$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "0.1203"
$ws.Range("E49").Value = "  -3.29%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "34.63"
$ws.Range("E50").Value = "  -3.79%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "0.05796"
$ws.Range("E51").Value = "  +0.44%  "

# Revert the cell style back to the workbook default (no explicit style
# index) now that the literal text values are safely stored, so the
# resulting cells match the original un-styled data cells.
$dataRange.Style = "Normal"
